$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing data (D..K) to (F..M)
$ws.Columns("D:E").Insert()

# Copy number formats from the (now-shifted) F:G columns into the new blank D:E columns
# so the new cells carry the same date / number style as their neighbours.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("G7").Value = 43190
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43008
$ws.Range("J7").Value = 42916
$ws.Range("K7").Value = 42825
$ws.Range("L7").Value = 42735
$ws.Range("M7").Value = 42643
$ws.Range("D8").Value = 540600
$ws.Range("E8").Value = 375200
$ws.Range("F8").Value = 367900
$ws.Range("G8").Value = 449900
$ws.Range("H8").Value = 1179000
$ws.Range("I8").Value = 247800
$ws.Range("J8").Value = 270600
$ws.Range("K8").Value = 326700
$ws.Range("L8").Value = 597800
$ws.Range("M8").Value = 499800
$ws.Range("D9").Value = 124700
$ws.Range("E9").Value = 103100
$ws.Range("F9").Value = 118600
$ws.Range("G9").Value = 149400
$ws.Range("H9").Value = 553700
$ws.Range("I9").Value = 369200
$ws.Range("J9").Value = 341300
$ws.Range("K9").Value = 346300
$ws.Range("L9").Value = 362900
$ws.Range("M9").Value = 330500
$ws.Range("D10").Value = 415900
$ws.Range("E10").Value = 272100
$ws.Range("F10").Value = 249300
$ws.Range("G10").Value = 300500
$ws.Range("H10").Value = 625300
$ws.Range("I10").Value = -121400
$ws.Range("J10").Value = -70700
$ws.Range("K10").Value = -19600
$ws.Range("L10").Value = 234900
$ws.Range("M10").Value = 169300
$ws.Range("D12").Value = 2600
$ws.Range("E12").Value = 3300
$ws.Range("F12").Value = 3700
$ws.Range("G12").Value = 2400
$ws.Range("H12").Value = 48100
$ws.Range("I12").Value = 4500
$ws.Range("J12").Value = 19700
$ws.Range("K12").Value = 9800
$ws.Range("L12").Value = 9500
$ws.Range("M12").Value = 400
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("D14").Value = -300
$ws.Range("E14").Value = 15400
$ws.Range("F14").Value = 42100
$ws.Range("G14").Value = 15600
$ws.Range("H14").Value = 140000
$ws.Range("I14").Value = 2000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 137000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("D15").Value = 130100
$ws.Range("E15").Value = 119600
$ws.Range("F15").Value = 119100
$ws.Range("G15").Value = 124700
$ws.Range("H15").Value = 412000
$ws.Range("I15").Value = 149200
$ws.Range("J15").Value = 133000
$ws.Range("K15").Value = 149100
$ws.Range("L15").Value = 156600
$ws.Range("M15").Value = 151700
$ws.Range("D17").Value = 300800
$ws.Range("E17").Value = 288400
$ws.Range("F17").Value = 328100
$ws.Range("G17").Value = -301400
$ws.Range("H17").Value = 1294000
$ws.Range("I17").Value = 285200
$ws.Range("J17").Value = 268500
$ws.Range("K17").Value = 405100
$ws.Range("L17").Value = 654500
$ws.Range("M17").Value = 582800
$ws.Range("D18").Value = 239800
$ws.Range("E18").Value = 86800
$ws.Range("F18").Value = 39800
$ws.Range("G18").Value = 751300
$ws.Range("H18").Value = -115000
$ws.Range("I18").Value = -37400
$ws.Range("J18").Value = 2100
$ws.Range("K18").Value = -78400
$ws.Range("L18").Value = -56700
$ws.Range("M18").Value = -83000
$ws.Range("D20").Value = -100900
$ws.Range("E20").Value = 152300
$ws.Range("F20").Value = 29000
$ws.Range("G20").Value = 46500
$ws.Range("H20").Value = 395000
$ws.Range("I20").Value = 65000
$ws.Range("J20").Value = 218400
$ws.Range("K20").Value = -18400
$ws.Range("L20").Value = -135800
$ws.Range("M20").Value = 245800
$ws.Range("D21").Value = 269000
$ws.Range("E21").Value = 358700
$ws.Range("F21").Value = 187800
$ws.Range("G21").Value = 922500
$ws.Range("H21").Value = 692100
$ws.Range("I21").Value = 129600
$ws.Range("J21").Value = 312100
$ws.Range("K21").Value = -1100
$ws.Range("L21").Value = -214500
$ws.Range("M21").Value = 314500
$ws.Range("D22").Value = 33200
$ws.Range("E22").Value = 35700
$ws.Range("F22").Value = 38400
$ws.Range("G22").Value = 38600
$ws.Range("H22").Value = 161400
$ws.Range("I22").Value = 38800
$ws.Range("J22").Value = 40700
$ws.Range("K22").Value = 41600
$ws.Range("L22").Value = 46900
$ws.Range("M22").Value = 47300
$ws.Range("D23").Value = 105700
$ws.Range("E23").Value = 203400
$ws.Range("F23").Value = 30300
$ws.Range("G23").Value = 759200
$ws.Range("H23").Value = 118600
$ws.Range("I23").Value = -11300
$ws.Range("J23").Value = 179800
$ws.Range("K23").Value = -138400
$ws.Range("L23").Value = -239400
$ws.Range("M23").Value = 115400
$ws.Range("D24").Value = -23700
$ws.Range("E24").Value = 56700
$ws.Range("F24").Value = -31100
$ws.Range("G24").Value = 213700
$ws.Range("H24").Value = -61200
$ws.Range("I24").Value = 10500
$ws.Range("J24").Value = 58000
$ws.Range("K24").Value = -47400
$ws.Range("L24").Value = 81800
$ws.Range("M24").Value = 52900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("D26").Value = 129400
$ws.Range("E26").Value = 146800
$ws.Range("F26").Value = 61400
$ws.Range("G26").Value = 545500
$ws.Range("H26").Value = 179700
$ws.Range("I26").Value = -21800
$ws.Range("J26").Value = 121800
$ws.Range("K26").Value = -91000
$ws.Range("L26").Value = -321200
$ws.Range("M26").Value = 62600
$ws.Range("D27").Value = 101900
$ws.Range("E27").Value = 125000
$ws.Range("F27").Value = 42000
$ws.Range("G27").Value = 527600
$ws.Range("H27").Value = 179700
$ws.Range("I27").Value = -21800
$ws.Range("J27").Value = 121800
$ws.Range("K27").Value = -91000
$ws.Range("L27").Value = -325600
$ws.Range("M27").Value = 60300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 201000
$ws.Range("I29").Value = -4600
$ws.Range("J29").Value = 47700
$ws.Range("K29").Value = 52000
$ws.Range("L29").Value = 19600
$ws.Range("M29").Value = -35000
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("D32").Value = 100900
$ws.Range("E32").Value = -152300
$ws.Range("F32").Value = -29000
$ws.Range("G32").Value = -46500
$ws.Range("H32").Value = -395000
$ws.Range("I32").Value = -65000
$ws.Range("J32").Value = -218400
$ws.Range("K32").Value = 18400
$ws.Range("L32").Value = 135800
$ws.Range("M32").Value = -245800
$ws.Range("D33").Value = 101900
$ws.Range("E33").Value = 125000
$ws.Range("F33").Value = 42000
$ws.Range("G33").Value = 527600
$ws.Range("H33").Value = 380700
$ws.Range("I33").Value = -26400
$ws.Range("J33").Value = 169500
$ws.Range("K33").Value = -39000
$ws.Range("L33").Value = -306000
$ws.Range("M33").Value = 25300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("D35").Value = 101900
$ws.Range("E35").Value = 125000
$ws.Range("F35").Value = 42000
$ws.Range("G35").Value = 527600
$ws.Range("H35").Value = 380700
$ws.Range("I35").Value = -26400
$ws.Range("J35").Value = 169500
$ws.Range("K35").Value = -39000
$ws.Range("L35").Value = -306000
$ws.Range("M35").Value = 25300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("G38").Value = 43190
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43008
$ws.Range("J38").Value = 42916
$ws.Range("K38").Value = 42825
$ws.Range("L38").Value = 42735
$ws.Range("M38").Value = 42643
$ws.Range("D41").Value = 17200
$ws.Range("E41").Value = 42700
$ws.Range("F41").Value = 54800
$ws.Range("G41").Value = 82500
$ws.Range("H41").Value = 509200
$ws.Range("I41").Value = 285700
$ws.Range("J41").Value = 299100
$ws.Range("K41").Value = 61300
$ws.Range("L41").Value = 46300
$ws.Range("M41").Value = 80200
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("D43").Value = 413000
$ws.Range("E43").Value = 197800
$ws.Range("F43").Value = 177200
$ws.Range("G43").Value = 221100
$ws.Range("H43").Value = 237200
$ws.Range("I43").Value = 377000
$ws.Range("J43").Value = 408800
$ws.Range("K43").Value = 440800
$ws.Range("L43").Value = 290100
$ws.Range("M43").Value = 244400
$ws.Range("D44").Value = 9700
$ws.Range("E44").Value = 9700
$ws.Range("F44").Value = 10500
$ws.Range("G44").Value = 10700
$ws.Range("H44").Value = 10700
$ws.Range("I44").Value = 63200
$ws.Range("J44").Value = 75000
$ws.Range("K44").Value = 69600
$ws.Range("L44").Value = 80800
$ws.Range("M44").Value = 62600
$ws.Range("D45").Value = 61800
$ws.Range("E45").Value = 65100
$ws.Range("F45").Value = 76600
$ws.Range("G45").Value = 92700
$ws.Range("H45").Value = 95300
$ws.Range("I45").Value = 79400
$ws.Range("J45").Value = 64200
$ws.Range("K45").Value = 44500
$ws.Range("L45").Value = 274400
$ws.Range("M45").Value = 127600
$ws.Range("D46").Value = 501700
$ws.Range("E46").Value = 315300
$ws.Range("F46").Value = 319100
$ws.Range("G46").Value = 406900
$ws.Range("H46").Value = 852500
$ws.Range("I46").Value = 805300
$ws.Range("J46").Value = 847100
$ws.Range("K46").Value = 616200
$ws.Range("L46").Value = 626100
$ws.Range("M46").Value = 514900
$ws.Range("D47").Value = 18700
$ws.Range("E47").Value = 19500
$ws.Range("F47").Value = 22300
$ws.Range("G47").Value = 20700
$ws.Range("H47").Value = 197900
$ws.Range("I47").Value = 190200
$ws.Range("J47").Value = 188600
$ws.Range("K47").Value = 197400
$ws.Range("L47").Value = 191000
$ws.Range("M47").Value = 257400
$ws.Range("D48").Value = 6942400
$ws.Range("E48").Value = 6768600
$ws.Range("F48").Value = 6541900
$ws.Range("G48").Value = 6621800
$ws.Range("H48").Value = 5789800
$ws.Range("I48").Value = 7799000
$ws.Range("J48").Value = 7794200
$ws.Range("K48").Value = 7822900
$ws.Range("L48").Value = 10902900
$ws.Range("M48").Value = 8414600
$ws.Range("D49").Value = 899600
$ws.Range("E49").Value = 901200
$ws.Range("F49").Value = 902800
$ws.Range("G49").Value = 923200
$ws.Range("H49").Value = "NA"
$ws.Range("I49").Value = "NA"
$ws.Range("J49").Value = "NA"
$ws.Range("K49").Value = "NA"
$ws.Range("L49").Value = "NA"
$ws.Range("M49").Value = "NA"
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("D52").Value = 229800
$ws.Range("E52").Value = 204400
$ws.Range("F52").Value = 421600
$ws.Range("G52").Value = 149600
$ws.Range("H52").Value = 91700
$ws.Range("I52").Value = 185200
$ws.Range("J52").Value = 195200
$ws.Range("K52").Value = 427300
$ws.Range("L52").Value = 2397900
$ws.Range("M52").Value = 378500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("D54").Value = 8592200
$ws.Range("E54").Value = 8209000
$ws.Range("F54").Value = 8207800
$ws.Range("G54").Value = 8122200
$ws.Range("H54").Value = 6931900
$ws.Range("I54").Value = 8979600
$ws.Range("J54").Value = 9025200
$ws.Range("K54").Value = 9063800
$ws.Range("L54").Value = 9179700
$ws.Range("M54").Value = 9565500
$ws.Range("D57").Value = 229800
$ws.Range("E57").Value = 263000
$ws.Range("F57").Value = 186400
$ws.Range("G57").Value = 193900
$ws.Range("H57").Value = 211200
$ws.Range("I57").Value = 303200
$ws.Range("J57").Value = 265100
$ws.Range("K57").Value = 270700
$ws.Range("L57").Value = 157100
$ws.Range("M57").Value = 197500
$ws.Range("D58").Value = 7000
$ws.Range("E58").Value = 7000
$ws.Range("F58").Value = 6900
$ws.Range("G58").Value = 6900
$ws.Range("H58").Value = 7100
$ws.Range("I58").Value = 11000
$ws.Range("J58").Value = 11400
$ws.Range("K58").Value = 11900
$ws.Range("L58").Value = 7900
$ws.Range("M58").Value = 358500
$ws.Range("D59").Value = 286200
$ws.Range("E59").Value = 263800
$ws.Range("F59").Value = 281700
$ws.Range("G59").Value = 236900
$ws.Range("H59").Value = 223400
$ws.Range("I59").Value = 546000
$ws.Range("J59").Value = 549200
$ws.Range("K59").Value = 710300
$ws.Range("L59").Value = 775000
$ws.Range("M59").Value = 514300
$ws.Range("D60").Value = 523000
$ws.Range("E60").Value = 533700
$ws.Range("F60").Value = 475000
$ws.Range("G60").Value = 437700
$ws.Range("H60").Value = 441700
$ws.Range("I60").Value = 860200
$ws.Range("J60").Value = 825700
$ws.Range("K60").Value = 992800
$ws.Range("L60").Value = 940000
$ws.Range("M60").Value = 1070200
$ws.Range("D61").Value = 2391500
$ws.Range("E61").Value = 2199600
$ws.Range("F61").Value = 2347600
$ws.Range("G61").Value = 2229800
$ws.Range("H61").Value = 2207400
$ws.Range("I61").Value = 2532300
$ws.Range("J61").Value = 2630100
$ws.Range("K61").Value = 2657300
$ws.Range("L61").Value = 2448400
$ws.Range("M61").Value = 2763800
$ws.Range("D62").Value = 595900
$ws.Range("E62").Value = 426900
$ws.Range("F62").Value = 346200
$ws.Range("G62").Value = 386900
$ws.Range("H62").Value = 383000
$ws.Range("I62").Value = 1515800
$ws.Range("J62").Value = 1481800
$ws.Range("K62").Value = 1507300
$ws.Range("L62").Value = 2467700
$ws.Range("M62").Value = 1441300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("D66").Value = 4262200
$ws.Range("E66").Value = 3898500
$ws.Range("F66").Value = 3899000
$ws.Range("G66").Value = 3778000
$ws.Range("H66").Value = 3032000
$ws.Range("I66").Value = 5047700
$ws.Range("J66").Value = 5079900
$ws.Range("K66").Value = 5300000
$ws.Range("L66").Value = 5381300
$ws.Range("M66").Value = 5418300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("D72").Value = 2071800
$ws.Range("E72").Value = 2003900
$ws.Range("F72").Value = 1940500
$ws.Range("G72").Value = 1940900
$ws.Range("H72").Value = 1455800
$ws.Range("I72").Value = 1825500
$ws.Range("J72").Value = 1852000
$ws.Range("K72").Value = 1682500
$ws.Range("L72").Value = 1727800
$ws.Range("M72").Value = 2033800
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("D76").Value = 4330000
$ws.Range("E76").Value = 4310600
$ws.Range("F76").Value = 4308800
$ws.Range("G76").Value = 4344200
$ws.Range("H76").Value = 3899900
$ws.Range("I76").Value = 3931800
$ws.Range("J76").Value = 3945300
$ws.Range("K76").Value = 3763800
$ws.Range("L76").Value = 3798400
$ws.Range("M76").Value = 4147200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("G80").Value = 43190
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43008
$ws.Range("J80").Value = 42916
$ws.Range("K80").Value = 42825
$ws.Range("L80").Value = 42735
$ws.Range("M80").Value = 42643
$ws.Range("D81").Value = 101900
$ws.Range("E81").Value = 125000
$ws.Range("F81").Value = 42000
$ws.Range("G81").Value = 527600
$ws.Range("H81").Value = 380700
$ws.Range("I81").Value = -26400
$ws.Range("J81").Value = 169500
$ws.Range("K81").Value = -39000
$ws.Range("L81").Value = -306000
$ws.Range("M81").Value = 25300
$ws.Range("D83").Value = 130100
$ws.Range("E83").Value = 119600
$ws.Range("F83").Value = 119100
$ws.Range("G83").Value = 124700
$ws.Range("H83").Value = 412000
$ws.Range("I83").Value = 148800
$ws.Range("J83").Value = 117100
$ws.Range("K83").Value = 148800
$ws.Range("L83").Value = 156600
$ws.Range("M83").Value = 151700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("D89").Value = 195600
$ws.Range("E89").Value = 239300
$ws.Range("F89").Value = 191600
$ws.Range("G89").Value = 259300
$ws.Range("H89").Value = 648700
$ws.Range("I89").Value = 181400
$ws.Range("J89").Value = 88700
$ws.Range("K89").Value = 211800
$ws.Range("L89").Value = 81000
$ws.Range("M89").Value = 164600
$ws.Range("D91").Value = -322300
$ws.Range("E91").Value = -297500
$ws.Range("F91").Value = -264200
$ws.Range("G91").Value = -232500
$ws.Range("H91").Value = -632800
$ws.Range("I91").Value = -149500
$ws.Range("J91").Value = -146000
$ws.Range("K91").Value = -216900
$ws.Range("L91").Value = -40800
$ws.Range("M91").Value = -243500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("D94").Value = -309800
$ws.Range("E94").Value = 54000
$ws.Range("F94").Value = -212500
$ws.Range("G94").Value = -426300
$ws.Range("H94").Value = -221900
$ws.Range("I94").Value = -84600
$ws.Range("J94").Value = 184200
$ws.Range("K94").Value = -87600
$ws.Range("L94").Value = 266100
$ws.Range("M94").Value = -70700
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("D100").Value = 88700
$ws.Range("E100").Value = -305500
$ws.Range("F100").Value = -6700
$ws.Range("G100").Value = -259700
$ws.Range("H100").Value = 36100
$ws.Range("I100").Value = -107100
$ws.Range("J100").Value = -35100
$ws.Range("K100").Value = -116700
$ws.Range("L100").Value = -368500
$ws.Range("M100").Value = -109600
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("D102").Value = -25500
$ws.Range("E102").Value = -12200
$ws.Range("F102").Value = -27600
$ws.Range("G102").Value = -426700
$ws.Range("H102").Value = 462900
$ws.Range("I102").Value = -10400
$ws.Range("J102").Value = 237800
$ws.Range("K102").Value = 7500
$ws.Range("L102").Value = -19800
$ws.Range("M102").Value = -17400
